# Update "想去人数" (column F) counts on multiple sheets to the values
# captured at the 456a3b4 gh-pages data refresh.
#
# NOTE: deliberately avoided "hashtable of arrays-of-pairs" here -
# PowerShell silently unrolls a single-element array literal
# (e.g. @(@(3, 2813)) collapses to @(3, 2813)), which would make a
# sheet with exactly one update desync row/value pairs. Using flat
# parallel arrays of (Sheet, Row, Value) triples sidesteps that.

$sheetCol = @(
    "展览","展览","展览","展览","展览","展览","展览","展览","展览","展览","展览",
    "展览","展览","展览","展览","展览","展览","展览","展览","展览","展览","展览",
    "演出","演出",
    "本地生活",
    "全部类型","全部类型","全部类型","全部类型","全部类型","全部类型","全部类型",
    "全部类型","全部类型","全部类型","全部类型","全部类型","全部类型","全部类型",
    "全部类型","全部类型","全部类型","全部类型","全部类型","全部类型","全部类型",
    "全部类型"
)

$rowCol = @(
    4,5,7,18,19,20,22,23,24,25,27,30,32,35,36,41,42,43,44,45,47,48,
    2,17,
    3,
    8,9,11,16,17,18,19,20,21,22,24,28,34,35,39,40,41,42,45,46,48,49
)

$valCol = @(
    429,8599,10807,77,14,414,1833,402,573,346,69,1211,4,1419,449,521,354,102,801,644,123,111,
    38,386,
    2813,
    429,8599,10807,77,14,414,1833,402,573,346,69,1211,1419,449,521,354,102,801,386,644,123,111
)

$wb = $excel.ActiveWorkbook

if ($sheetCol.Count -ne $rowCol.Count -or $sheetCol.Count -ne $valCol.Count) {
    throw "Update arrays are out of sync: $($sheetCol.Count) / $($rowCol.Count) / $($valCol.Count)"
}

for ($i = 0; $i -lt $sheetCol.Count; $i++) {
    $ws = $wb.Worksheets.Item($sheetCol[$i])
    $ws.Cells.Item($rowCol[$i], 6).Value = $valCol[$i]
}
